$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.811.15"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.70"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5021"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2564"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06380"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.59"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07692"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.641.09"
$ws.Range("E12").Value = "  -0.81%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.240"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").Value = "1.858.10"
$ws.Range("E14").Value = "  -1.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5419"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").Value = "0.0₅7896"
$ws.Range("E16").Value = "  -1.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.49"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "25.827.57"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("E20").Value = "  -3.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.327"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.915"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.953"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.87%  "
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.929"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +11.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.50"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.38%  "
$ws.Range("E27").Value = "  -2.45%  "
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.697"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.238"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04989"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.257"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.173"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("E34").Value = "  -1.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.363"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.48%  "
$ws.Range("D36").Value = "1.170.05"
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("E37").Value = "  -4.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.612"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.95%  "
$ws.Range("E39").Value = "  -2.02%  "
$ws.Range("E40").Value = "  -2.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.556"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.678"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8062"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.29"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("D46").Value = "1.770.22"
$ws.Range("D47").Value = "0.0₈115"
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4512"
$ws.Range("D48").ClearFormats()
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.62"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05078"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.77%  "
